$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 32370
$ws.Cells.Item(3, 1).Value = 32462
$ws.Cells.Item(4, 1).Value = 32554
$ws.Cells.Item(5, 1).Value = 32643
$ws.Cells.Item(6, 1).Value = 32735
$ws.Cells.Item(7, 1).Value = 32827
$ws.Cells.Item(8, 1).Value = 32919
$ws.Cells.Item(9, 1).Value = 33008
$ws.Cells.Item(10, 1).Value = 33100
$ws.Cells.Item(11, 1).Value = 33192
$ws.Cells.Item(12, 1).Value = 33284
$ws.Cells.Item(13, 1).Value = 33373
$ws.Cells.Item(14, 1).Value = 33465
$ws.Cells.Item(15, 1).Value = 33557
$ws.Cells.Item(16, 1).Value = 33649
$ws.Cells.Item(17, 1).Value = 33739
$ws.Cells.Item(18, 1).Value = 33831
$ws.Cells.Item(19, 1).Value = 33923
$ws.Cells.Item(20, 1).Value = 34015
$ws.Cells.Item(21, 1).Value = 34104
$ws.Cells.Item(22, 1).Value = 34196
$ws.Cells.Item(23, 1).Value = 34288
$ws.Cells.Item(24, 1).Value = 34380
$ws.Cells.Item(25, 1).Value = 34469
$ws.Cells.Item(26, 1).Value = 34561
$ws.Cells.Item(27, 1).Value = 34653
$ws.Cells.Item(28, 1).Value = 34745
$ws.Cells.Item(29, 1).Value = 34834
$ws.Cells.Item(30, 1).Value = 34926
$ws.Cells.Item(31, 1).Value = 35018
$ws.Cells.Item(32, 1).Value = 35110
$ws.Cells.Item(33, 1).Value = 35200
$ws.Cells.Item(34, 1).Value = 35292
$ws.Cells.Item(35, 1).Value = 35384
$ws.Cells.Item(36, 1).Value = 35476
$ws.Cells.Item(37, 1).Value = 35565
$ws.Cells.Item(38, 1).Value = 35657
$ws.Cells.Item(39, 1).Value = 35749
$ws.Cells.Item(40, 1).Value = 35841
$ws.Cells.Item(41, 1).Value = 35930
$ws.Cells.Item(42, 1).Value = 36022
$ws.Cells.Item(43, 1).Value = 36114
$ws.Cells.Item(44, 1).Value = 36206
$ws.Cells.Item(45, 1).Value = 36295
$ws.Cells.Item(46, 1).Value = 36387
$ws.Cells.Item(47, 1).Value = 36479
$ws.Cells.Item(48, 1).Value = 36571
$ws.Cells.Item(49, 1).Value = 36661
$ws.Cells.Item(50, 1).Value = 36753
$ws.Cells.Item(51, 1).Value = 36845
$ws.Cells.Item(52, 1).Value = 36937
$ws.Cells.Item(53, 1).Value = 37026
$ws.Cells.Item(54, 1).Value = 37118
$ws.Cells.Item(55, 1).Value = 37210
$ws.Cells.Item(56, 1).Value = 37302
$ws.Cells.Item(57, 1).Value = 37391
$ws.Cells.Item(58, 1).Value = 37483
$ws.Cells.Item(59, 1).Value = 37575
$ws.Cells.Item(60, 1).Value = 37667
$ws.Cells.Item(61, 1).Value = 37756
$ws.Cells.Item(62, 1).Value = 37848
$ws.Cells.Item(63, 1).Value = 37940
$ws.Cells.Item(64, 1).Value = 38032
$ws.Cells.Item(65, 1).Value = 38122
$ws.Cells.Item(66, 1).Value = 38214
$ws.Cells.Item(67, 1).Value = 38306
$ws.Cells.Item(68, 1).Value = 38398
$ws.Cells.Item(69, 1).Value = 38487
$ws.Cells.Item(70, 1).Value = 38579
$ws.Cells.Item(71, 1).Value = 38671
$ws.Cells.Item(72, 1).Value = 38763
$ws.Cells.Item(73, 1).Value = 38852
$ws.Cells.Item(74, 1).Value = 38944
$ws.Cells.Item(75, 1).Value = 39036
$ws.Cells.Item(76, 1).Value = 39128
$ws.Cells.Item(77, 1).Value = 39217
$ws.Cells.Item(78, 1).Value = 39309
$ws.Cells.Item(79, 1).Value = 39401
$ws.Cells.Item(80, 1).Value = 39493
$ws.Cells.Item(81, 1).Value = 39583
$ws.Cells.Item(82, 1).Value = 39675
$ws.Cells.Item(83, 1).Value = 39767
$ws.Cells.Item(84, 1).Value = 39859
$ws.Cells.Item(85, 1).Value = 39948
$ws.Cells.Item(86, 1).Value = 40040
$ws.Cells.Item(87, 1).Value = 40132
$ws.Cells.Item(88, 1).Value = 40224
$ws.Cells.Item(89, 1).Value = 40313
$ws.Cells.Item(90, 1).Value = 40405
$ws.Cells.Item(91, 1).Value = 40497
$ws.Cells.Item(92, 1).Value = 40589
$ws.Cells.Item(93, 1).Value = 40678
$ws.Cells.Item(94, 1).Value = 40770
$ws.Cells.Item(95, 1).Value = 40862
$ws.Cells.Item(96, 1).Value = 40954
$ws.Cells.Item(97, 1).Value = 41044
$ws.Cells.Item(98, 1).Value = 41136
$ws.Cells.Item(99, 1).Value = 41228
$ws.Cells.Item(100, 1).Value = 41320
$ws.Cells.Item(101, 1).Value = 41409
$ws.Cells.Item(102, 1).Value = 41501
$ws.Cells.Item(103, 1).Value = 41593
$ws.Cells.Item(104, 1).Value = 41685
$ws.Cells.Item(105, 1).Value = 41774
$ws.Cells.Item(106, 1).Value = 41866
$ws.Cells.Item(107, 1).Value = 41958
$ws.Cells.Item(108, 1).Value = 42050
$ws.Cells.Item(109, 1).Value = 42139
$ws.Cells.Item(110, 1).Value = 42231
$ws.Cells.Item(111, 1).Value = 42323
$ws.Cells.Item(112, 1).Value = 42415
$ws.Cells.Item(113, 1).Value = 42505
$ws.Cells.Item(114, 1).Value = 42597
$ws.Cells.Item(115, 1).Value = 42689
$ws.Cells.Item(116, 1).Value = 42781
$ws.Cells.Item(117, 1).Value = 42870
$ws.Cells.Item(118, 1).Value = 42962
$ws.Cells.Item(119, 1).Value = 43054
$ws.Cells.Item(120, 1).Value = 43146
$ws.Cells.Item(121, 1).Value = 43235
$ws.Cells.Item(122, 1).Value = 43327
$ws.Cells.Item(123, 1).Value = 43419
$ws.Cells.Item(124, 1).Value = 43511
$ws.Cells.Item(125, 1).Value = 43600
$ws.Cells.Item(126, 1).Value = 43692
$ws.Cells.Item(127, 1).Value = 43784
$ws.Cells.Item(128, 1).Value = 43876
$ws.Cells.Item(129, 1).Value = 43966
$ws.Cells.Item(130, 1).Value = 44058
$ws.Cells.Item(131, 1).Value = 44150
$ws.Cells.Item(132, 1).Value = 44242
$ws.Cells.Item(133, 1).Value = 44331
$ws.Cells.Item(134, 1).Value = 44423
$ws.Cells.Item(135, 1).Value = 44515
$ws.Cells.Item(136, 1).Value = 44607
$ws.Cells.Item(137, 1).Value = 44696
$ws.Cells.Item(138, 1).Value = 44788
$ws.Cells.Item(139, 1).Value = 44880
$ws.Cells.Item(140, 1).Value = 44972
$ws.Cells.Item(141, 1).Value = 45061
$ws.Cells.Item(142, 1).Value = 45153
$ws.Cells.Item(143, 1).Value = 45245
$ws.Cells.Item(144, 1).Value = 45337
$ws.Cells.Item(145, 1).Value = 45427
$ws.Cells.Item(146, 1).Value = 45519
$ws.Cells.Item(147, 1).Value = 45611
$ws.Cells.Item(148, 1).Value = 45703
$ws.Cells.Item(149, 1).Value = 45792
$ws.Cells.Item(150, 1).Value = 45884
